$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("A2").Value = 112390694
$ws.Range("B2").Value = 90806
$ws.Range("Q2").Value = 508632
$ws.Range("R2").Value = 6784369

# Row 3 update
$ws.Range("B3").Value = 88180

# Row 4 updates
$ws.Range("A4").Value = 112390653
$ws.Range("B4").Value = 90806
$ws.Range("Q4").Value = 508942
$ws.Range("R4").Value = 6784419
